$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C:F across rows 3-5 (cyclic shift: row5 -> row3, row3 -> row4, row4 -> row5)
$ws.Range("C3").Value = "5"
$ws.Range("D3").Value = "8"
$ws.Range("E3").Value = "0"
$ws.Range("F3").Value = "0"

$ws.Range("C4").Value = "26"
$ws.Range("D4").Value = "10"
$ws.Range("E4").Value = "1"
$ws.Range("F4").Value = "3"

$ws.Range("C5").Value = "11"
$ws.Range("D5").Value = "15"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "0"
